# Update Andres Gil's "ULTIMO_PAGO" (last payment) date in D2 to a new,
# differently-formatted date string ("2023/12/25" instead of "28/12/2023").
# Force text storage (not an auto-converted date serial) by temporarily
# switching the cell to a text number format, matching the original file
# where dates are kept as plain text (shared strings), then restore the
# default "General" format so the cell's style index is unaffected.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "2023/12/25"
$ws.Range("D2").NumberFormat = "General"

# Reflect the new active cell/selection shown in the sheet view (moved
# from D11 to D3).
$ws.Range("D3").Select()
